# Add 2022-Q3 data
#
# 1. Insert a brand-new worksheet "2022-Q3" right before the existing
#    "2022-Q2" sheet and populate it with the latest quarter's fund-holding
#    detail table (same layout/headers as the other quarterly sheets).
# 2. Insert a new top data-row in the "总计" (summary) sheet for 2022-Q3 and
#    shift the existing history down by one row, renumbering the index
#    column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: write one header cell (bold / centered / boxed), matching the
# style already used for the header row + index column on the other
# quarterly sheets.
# ---------------------------------------------------------------------
function Set-BoxedCell($ws, $row, $col, $val) {
    $c = $ws.Cells.Item($row, $col)
    $c.Value = $val
    $c.Font.Bold = $true
    $c.HorizontalAlignment = -4108   # xlCenter
    $c.VerticalAlignment = -4160     # xlTop
    $c.Borders.Item(7).LineStyle = 1   # xlEdgeLeft
    $c.Borders.Item(8).LineStyle = 1   # xlEdgeTop
    $c.Borders.Item(9).LineStyle = 1   # xlEdgeBottom
    $c.Borders.Item(10).LineStyle = 1  # xlEdgeRight
}

# ---------------------------------------------------------------------
# 1. New "2022-Q3" worksheet, inserted before "2022-Q2"
# ---------------------------------------------------------------------
$q2Sheet = $wb.Worksheets.Item("2022-Q2")
$q3Sheet = $wb.Worksheets.Add($q2Sheet)
$q3Sheet.Name = "2022-Q3"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
$col = 2
foreach ($h in $headers) {
    Set-BoxedCell $q3Sheet 1 $col $h
    $col = $col + 1
}

# Row 2 - 合煦智远嘉选混合A
Set-BoxedCell $q3Sheet 2 1 0
$q3Sheet.Cells.Item(2, 2).Value = "'006323"
$q3Sheet.Cells.Item(2, 3).Value = "合煦智远嘉选混合A"
$q3Sheet.Cells.Item(2, 4).Value = "'0.78"
$q3Sheet.Cells.Item(2, 5).Value = "'73.07"
$q3Sheet.Cells.Item(2, 6).Value = "'4.11"
$q3Sheet.Cells.Item(2, 7).Value = "'0.0321"
$q3Sheet.Cells.Item(2, 8).Value = 4

# Row 3 - 合煦智远嘉选混合C
Set-BoxedCell $q3Sheet 3 1 1
$q3Sheet.Cells.Item(3, 2).Value = "'006324"
$q3Sheet.Cells.Item(3, 3).Value = "合煦智远嘉选混合C"
$q3Sheet.Cells.Item(3, 4).Value = "'0.14"
$q3Sheet.Cells.Item(3, 5).Value = "'73.07"
$q3Sheet.Cells.Item(3, 6).Value = "'4.11"
$q3Sheet.Cells.Item(3, 7).Value = "'0.0058"
$q3Sheet.Cells.Item(3, 8).Value = 4

# ---------------------------------------------------------------------
# 2. "总计" summary sheet - insert a new row for 2022-Q3 at the top of the
#    data table and shift the rest down, renumbering the index column.
# ---------------------------------------------------------------------
$zj = $wb.Worksheets.Item("总计")
$zj.Rows.Item(2).Insert()
$zj.Range("B2:D2").ClearFormats()

# Copy the index-column box style down onto the newly inserted A2 cell.
$zj.Range("A3").Copy()
$zj.Range("A2").PasteSpecial(-4122)   # xlPasteFormats

$zj.Range("A2").Value = 0
$zj.Range("B2").Value = "2022-Q3"
$zj.Range("C2").Value = 2
$zj.Range("D2").Value = 0.04

# Renumber the index column for the rows that shifted down (were 0..6 in
# rows 2..8, now need to be 1..7 in rows 3..9).
for ($r = 3; $r -le 9; $r++) {
    $zj.Cells.Item($r, 1).Value = $r - 2
}

# Restore the originally-selected tab (last sheet, "2020-Q4").
$q4Sheet = $wb.Worksheets.Item("2020-Q4")
$q4Sheet.Activate()
$q4Sheet.Range("A1").Select()
